# feat: add 2022-Q4 data
#
# 1. Duplicate the existing "2022-Q3" worksheet (this keeps all of its
#    sheet-level formatting/page setup) to create the new "2022-Q4" sheet
#    positioned right after "总计", then overwrite its cells with the
#    2022-Q4 fund-holding detail rows (it only needs 4 data rows, so the
#    leftover 6th/7th template rows are cleared).
# 2. On the "总计" (summary) sheet, shift every quarter's B:D values down
#    one row and write the new 2022-Q4 summary values into row 2. Column A
#    (the 0-based row index) is left untouched, and a new row is appended
#    at the bottom carrying the values that used to belong to the last
#    quarter (2020-Q4).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2022-Q3" into a new "2022-Q4" sheet placed right
# before it (i.e. right after "总计").
# ---------------------------------------------------------------------
$q3ForCopy = $wb.Worksheets.Item("2022-Q3")
$q3ForCopy.Copy($q3ForCopy)

# NOTE: worksheet objects captured before a sheet-collection mutation can
# become stale, so everything we still need is re-fetched by name.
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"
$newSheet = $wb.Worksheets.Item("2022-Q4")

# The template ("2022-Q3") has 6 data rows; 2022-Q4 only needs 4, so drop
# the leftover rows 6 and 7.
$newSheet.Range("A6:H7").Clear()

# Columns that must stay text even though their contents look numeric
# (fund code / size / position values keep leading zeros & fixed decimals).
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$fundRows = @(
    @{ Row = 2; A = 0; Code = "004702"; Name = "南方金融主题灵活配置混合A"; Size = "12.97"; Position = "92.71"; Ratio = "4.95"; Value = "0.6420"; Rank = 5 },
    @{ Row = 3; A = 1; Code = "013500"; Name = "南方金融主题灵活配置混合C"; Size = "4.80";  Position = "92.71"; Ratio = "4.95"; Value = "0.2376"; Rank = 5 },
    @{ Row = 4; A = 2; Code = "398041"; Name = "中海量化策略混合";          Size = "2.41";  Position = "91.05"; Ratio = "5.98"; Value = "0.1441"; Rank = 5 },
    @{ Row = 5; A = 3; Code = "515760"; Name = "华夏中证浙江国资创新发展ETF"; Size = "2.14"; Position = "99.05"; Ratio = "4.53"; Value = "0.0969"; Rank = 7 }
)

foreach ($r in $fundRows) {
    $row = $r.Row
    $newSheet.Range("A$row").Value = $r.A
    $newSheet.Range("B$row").Value = $r.Code
    $newSheet.Range("C$row").Value = $r.Name
    $newSheet.Range("D$row").Value = $r.Size
    $newSheet.Range("E$row").Value = $r.Position
    $newSheet.Range("F$row").Value = $r.Ratio
    $newSheet.Range("G$row").Value = $r.Value
    $newSheet.Range("H$row").Value = $r.Rank
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet — push every quarter's B:D
# values down one row, write the new 2022-Q4 row at the top, and append
# the row that drops off the bottom (2020-Q4).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

$summaryRows = @(
    @{ Row = 2; Quarter = "2022-Q4"; Count = 4;  Aum = 1.12 },
    @{ Row = 3; Quarter = "2022-Q3"; Count = 6;  Aum = 0.71 },
    @{ Row = 4; Quarter = "2022-Q2"; Count = 3;  Aum = 0.15 },
    @{ Row = 5; Quarter = "2021-Q4"; Count = 10; Aum = 0.39 },
    @{ Row = 6; Quarter = "2021-Q3"; Count = 14; Aum = 0.75 },
    @{ Row = 7; Quarter = "2021-Q2"; Count = 6;  Aum = 0.68 },
    @{ Row = 8; Quarter = "2021-Q1"; Count = 5;  Aum = 0.45 },
    @{ Row = 9; Quarter = "2020-Q4"; Count = 5;  Aum = 0.59 }
)

foreach ($r in $summaryRows) {
    $row = $r.Row
    $totalSheet.Range("B$row").Value = $r.Quarter
    $totalSheet.Range("C$row").Value = $r.Count
    $totalSheet.Range("D$row").Value = $r.Aum
}

# Column A keeps its original 0-based sequence; extend it to the new row.
$totalSheet.Range("A9").Value = 7

# Restore "总计" as the active sheet/tab (duplicating a sheet makes the
# new one active by default).
$totalSheet.Activate()
